# Apply missing-data imputations / re-blankings and remove two rows
# (RM 232 and SC 92), causing the remaining rows to shift up.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Cell-level imputations / blankings (row numbers as in the original sheet) ---
$ws.Range("C3").Value = 11.2      # RM 8  : impute B value
$ws.Range("F4").Value = ""        # RM 9  : blank out F value
$ws.Range("C5").Value = ""        # RM 14 : blank out B value
$ws.Range("F9").Value = 17.26     # RM 42 : impute F value
$ws.Range("F10").Value = 16.43    # RM 52a: impute F value
$ws.Range("F17").Value = ""       # RM 116: blank out F value
$ws.Range("F18").Value = ""       # RM 120: blank out F value
$ws.Range("C21").Value = 12.7     # RM 135: impute B value
$ws.Range("C23").Value = ""       # RM 140: blank out B value

# --- Remove row 26 (RM 232) entirely ---
$ws.Rows(26).Delete()

# After the row-26 deletion, the row that used to be 28 (SC 92) is now row 27.
$ws.Rows(27).Delete()

# --- Impute the B value for SC 193, which is now row 32 ---
$ws.Range("C32").Value = 10.5
